$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.991.57"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "2.197.34"
$ws.Range("E3").Value = "  -2.22%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "295.43"
$ws.Range("E5").Value = "  -4.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "89.60"
$ws.Range("E6").Value = "  -5.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.563"
$ws.Range("E7").Value = "  -1.59%  "
$ws.Range("E9").Value = "  -7.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.44"
$ws.Range("E10").Value = "  -7.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0775"
$ws.Range("E11").Value = "  -4.65%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.88"
$ws.Range("E12").Value = "  -4.96%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.103"
$ws.Range("E13").Value = "  -0.96%  "
$ws.Range("D14").Value = "2.531.70"
$ws.Range("E14").Value = "  -2.27%  "
$ws.Range("D15").Value = "2.257.44"
$ws.Range("E15").Value = "  -4.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.28"
$ws.Range("E16").Value = "  -3.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.780"
$ws.Range("E17").Value = "  -7.57%  "
$ws.Range("D18").Value = "43.667.66"
$ws.Range("E18").Value = "  -1.05%  "
$ws.Range("E19").Value = "  -7.79%  "
$ws.Range("E20").Value = "  -9.12%  "
$ws.Range("E21").Value = "  -11.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "63.35"
$ws.Range("E22").Value = "  -4.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.02"
$ws.Range("E23").Value = "  -2.11%  "
$ws.Range("E24").Value = "  -13.90%  "
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("E26").Value = "  -8.77%  "
$ws.Range("E27").Value = "  +0.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.97"
$ws.Range("E28").Value = "  -4.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.28"
$ws.Range("E29").Value = "  -6.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.31"
$ws.Range("E30").Value = "  -4.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "149.16"
$ws.Range("E31").Value = "  -3.49%  "
$ws.Range("E32").Value = "  -10.99%  "
$ws.Range("E33").Value = "  -4.85%  "
$ws.Range("E34").Value = "  -7.91%  "
$ws.Range("E35").Value = "  -3.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.84"
$ws.Range("E36").Value = "  -9.87%  "
$ws.Range("E37").Value = "  -6.13%  "
$ws.Range("E38").Value = "  -9.00%  "
$ws.Range("E39").Value = "  -5.96%  "
$ws.Range("E40").Value = "  -8.28%  "
$ws.Range("E41").Value = "  -12.19%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.01"
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.05"
$ws.Range("E43").Value = "  -11.55%  "
$ws.Range("D44").Value = "1.797.34"
$ws.Range("E44").Value = "  +2.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.70"
$ws.Range("E45").Value = "  +5.74%  "
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.176"
$ws.Range("E46").Value = "  -9.88%  "
$ws.Range("E47").Value = "  +13.56%  "
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "72.65"
$ws.Range("E48").Value = "  -10.36%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "92.26"
$ws.Range("E49").Value = "  -7.92%  "
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "65.25"
$ws.Range("E50").Value = "  -8.42%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.414.98"
$ws.Range("E51").Value = "  -2.18%  "
